# WAT API automation test scripts
# Adds a new validation row (WAT-541) to the WoS_AuthorTransformation sheet,
# and backfills the STORE column (K) for two existing rows (WAT-429 / WAT-430).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Backfill STORE (col K) on existing rows -----------------------------
# Order matters: new shared-string entries must land in this exact sequence
# (lastName||firstName, status=200||..., orcid, Verify..., WAT-541).
$ws.Range("K39").Value = "lastName||firstName"

$ws.Range("J43").Value = "status=200||hits[0].primaryName=(WAT-430_lastName)||hits[0].primaryName=(WAT-430_firstName)"

$ws.Range("K34").Value = "orcid"

$ws.Range("B43").Value = "Verify that author cluster details primary name should match with author metadata last name and first name"

$ws.Range("A43").Value = "WAT-541"

# --- 2) Give the new row 43 the same look & feel as the rows above it ------
$ws.Range("A42").Copy()
$ws.Range("A43").PasteSpecial(-4122)

$ws.Range("B42").Copy()
$ws.Range("B43").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("C43").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D43").PasteSpecial(-4122)

$ws.Range("L2").Copy()
$ws.Range("E43").PasteSpecial(-4122)

$ws.Range("F42").Copy()
$ws.Range("F43").PasteSpecial(-4122)

$ws.Range("K34").Copy()
$ws.Range("G43").PasteSpecial(-4122)

$ws.Range("H42").Copy()
$ws.Range("H43").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("I43").PasteSpecial(-4122)

$ws.Range("J42").Copy()
$ws.Range("J43").PasteSpecial(-4122)

$ws.Range("K34").Copy()
$ws.Range("K43").PasteSpecial(-4122)

$ws.Rows(43).RowHeight = 30

# --- 3) Fill in the rest of row 43's data -----------------------------------
$ws.Range("C43").Value = "1PRECOMMEND"
$ws.Range("D43").Value = "/recommend/search/author/clusters/0000-0002-6423-7213"
$ws.Range("E43").Value = "GET"
$ws.Range("I43").Value = "WAT-430"

# --- 4) L42 and L43 both read PASS, with the "no explicit style" look ------
# (clear L42's format first so it matches the un-styled look of L41/L43)
$ws.Range("L41").Copy()
$ws.Range("L42").PasteSpecial(-4122)
$ws.Range("L42").Value = "PASS"

$ws.Range("L41").Copy()
$ws.Range("L43").PasteSpecial(-4122)
$ws.Range("L43").Value = "PASS"

# --- 5) Selection follows the newly-added row -------------------------------
$null = $ws.Range("A43").Select()
